# Auto-generated from XML diff: update NATMI TPM-derived metrics for Dcn-Egfr pairs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.233789666666667
$ws.Range("H2").Value = 3.701369
$ws.Range("I2").Value = 0.0001664233864291757
$ws.Range("J2").Value = 0.0001664233864291757
$ws.Range("M2").Value = 0.4102596666666667
$ws.Range("N2").Value = 1.230779
$ws.Range("O2").Value = 0.003499619873322347
$ws.Range("P2").Value = 0.003499619873322347
$ws.Range("Q2").Value = 0.5061741373834445
$ws.Range("R2").Value = 4.555567236451
$ws.Range("S2").Value = [double]"5.824185905331481E-07"
$ws.Range("T2").Value = [double]"5.82418590533148E-07"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.233789666666667
$ws.Range("H3").Value = 3.701369
$ws.Range("I3").Value = 0.0001664233864291757
$ws.Range("J3").Value = 0.0001664233864291757
$ws.Range("O3").Value = 0.8692174743460166
$ws.Range("P3").Value = 0.8692174743460165
$ws.Range("Q3").Value = 125.7209128996123
$ws.Range("R3").Value = 1131.488216096511
$ws.Range("S3").Value = 0.0001446581156240792
$ws.Range("T3").Value = 0.0001446581156240792
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.233789666666667
$ws.Range("H4").Value = 3.701369
$ws.Range("I4").Value = 0.0001664233864291757
$ws.Range("J4").Value = 0.0001664233864291757
$ws.Range("N4").Value = 44.764041
$ws.Range("O4").Value = 0.1272829057806611
$ws.Range("P4").Value = 0.1272829057806611
$ws.Range("Q4").Value = 18.40980374134767
$ws.Range("R4").Value = 165.688233672129
$ws.Range("S4").Value = [double]"2.118285221456333E-05"
$ws.Range("T4").Value = [double]"2.118285221456333E-05"
$ws.Range("I5").Value = 0.9827534361704352
$ws.Range("J5").Value = 0.9827534361704352
$ws.Range("M5").Value = 0.4102596666666667
$ws.Range("N5").Value = 1.230779
$ws.Range("O5").Value = 0.003499619873322347
$ws.Range("P5").Value = 0.003499619873322347
$ws.Range("Q5").Value = 2989.029267385337
$ws.Range("R5").Value = 26901.26340646804
$ws.Range("S5").Value = 0.00343926345579788
$ws.Range("T5").Value = 0.00343926345579788
$ws.Range("I6").Value = 0.9827534361704352
$ws.Range("J6").Value = 0.9827534361704352
$ws.Range("O6").Value = 0.8692174743460166
$ws.Range("P6").Value = 0.8692174743460165
$ws.Range("S6").Value = 0.8542264596929349
$ws.Range("T6").Value = 0.8542264596929348
$ws.Range("I7").Value = 0.9827534361704352
$ws.Range("J7").Value = 0.9827534361704352
$ws.Range("N7").Value = 44.764041
$ws.Range("O7").Value = 0.1272829057806611
$ws.Range("P7").Value = 0.1272829057806611
$ws.Range("R7").Value = 978412.2560418522
$ws.Range("S7").Value = 0.1250877130217025
$ws.Range("T7").Value = 0.1250877130217025
$ws.Range("I8").Value = 0.01708014044313564
$ws.Range("J8").Value = 0.01708014044313564
$ws.Range("M8").Value = 0.4102596666666667
$ws.Range("N8").Value = 1.230779
$ws.Range("O8").Value = 0.003499619873322347
$ws.Range("P8").Value = 0.003499619873322347
$ws.Range("Q8").Value = 51.94898109390123
$ws.Range("R8").Value = 467.540829845111
$ws.Range("S8").Value = [double]"5.977399893393426E-05"
$ws.Range("T8").Value = [double]"5.977399893393425E-05"
$ws.Range("I9").Value = 0.01708014044313564
$ws.Range("J9").Value = 0.01708014044313564
$ws.Range("O9").Value = 0.8692174743460166
$ws.Range("P9").Value = 0.8692174743460165
$ws.Range("S9").Value = 0.01484635653745762
$ws.Range("T9").Value = 0.01484635653745761
$ws.Range("I10").Value = 0.01708014044313564
$ws.Range("J10").Value = 0.01708014044313564
$ws.Range("N10").Value = 44.764041
$ws.Range("O10").Value = 0.1272829057806611
$ws.Range("P10").Value = 0.1272829057806611
$ws.Range("Q10").Value = 1889.410137478474
$ws.Range("S10").Value = 0.002174009906744094
$ws.Range("T10").Value = 0.002174009906744094

Write-Host "Updated 87 cells"
